# This script reproduces the authoring edit captured in the commit:
# "case18_1, case18_2, case18_3. 2020 and 2025. RES added."
#
# The core content change is that the "RES installed" sheet's installed
# capacity column (C2:C6) is populated with non-zero values (previously
# all zeros). Every other change in the workbook (Main!B7 total, and the
# volatile VLOOKUP/RANDBETWEEN-driven values on the "Pg, Winter, S#" and
# "Pg, Summer, S#" sheets) is simply the natural recalculation cascade
# that results from that single input edit, so we just need to set the
# inputs and force Excel to recalculate.
#
# Additionally, the active/selected sheet tab moves from "RES installed"
# back to "Main".

$wb = $excel.ActiveWorkbook

# --- Update the RES installed capacities (the actual authored edit) ---
$resWs = $wb.Worksheets.Item("RES installed")
$resWs.Range("C2").Value = 2
$resWs.Range("C3").Value = 5
$resWs.Range("C4").Value = 1
$resWs.Range("C5").Value = 1
$resWs.Range("C6").Value = 1

# --- Force a full workbook recalculation so every dependent formula ---
# --- (Main!B7 SUM, and all the VLOOKUP(...)*AVERAGE(...)*RANDBETWEEN() ---
# --- cells on the Pg sheets) picks up the new, non-zero RES capacity. ---
$excel.CalculateFullRebuild()
$excel.Calculate()

# --- Restore the originally-intended active tab: "Main" ---
$mainWs = $wb.Worksheets.Item("Main")
$mainWs.Activate()
$mainWs.Range("B5").Select()
